$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (rows 4-11 shift down to 5-12)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new data record
$ws.Cells.Item(4,1).Value = 11
$ws.Cells.Item(4,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4,3).Value = "Bíobío"
$ws.Cells.Item(4,4).Value = 44901
$ws.Cells.Item(4,5).Value = 8
$ws.Cells.Item(4,6).Value = "Fruta"
$ws.Cells.Item(4,7).Value = 100103
$ws.Cells.Item(4,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(4,9).Value = 100103003
$ws.Cells.Item(4,10).Value = "Damasco"
$ws.Cells.Item(4,11).Value = "Castle Brite"
$ws.Cells.Item(4,12).Value = "Primera"
$ws.Cells.Item(4,13).Value = 100
$ws.Cells.Item(4,14).Value = 15000
$ws.Cells.Item(4,15).Value = 16000
$ws.Cells.Item(4,16).Value = 15500
$ws.Cells.Item(4,17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(4,18).Value = "Región de O'Higgins"
$ws.Cells.Item(4,19).Value = 1550
$ws.Cells.Item(4,20).Value = 10
